{"js": "// Remove the two blank spacer paragraphs and the \"\u00a9 2020 ...\" footer\n// paragraph that used to immediately follow the\n// \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph that contains the LOQ4086 requirement line.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4086:\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4086' anchor paragraph.\");\n}\n\n// The three paragraphs right after the anchor are: an empty spacer\n// paragraph, an empty page-break spacer paragraph, and the\n// \"\u00a9 2020 ... Contact: ... Jekyll ...\" copyright/footer paragraph.\n// Delete exactly those three, identified by their known content so the\n// script stays correct even if surrounding spacing shifts slightly.\nconst candidates = [];\nfor (let i = anchorIndex + 1; i < items.length && candidates.length < 3; i++) {\n  const text = items[i].text.trim();\n  const isBlankSpacer = text === \"\";\n  const isCopyrightFooter = text.indexOf(\"Contact:\") !== -1 && text.indexOf(\"Jekyll\") !== -1;\n  if (isBlankSpacer || isCopyrightFooter) {\n    candidates.push(items[i]);\n  } else {\n    break;\n  }\n}\n\nif (candidates.length !== 3) {\n  throw new Error(\"Expected exactly 3 paragraphs to remove after the LOQ4086 anchor, found \" + candidates.length);\n}\n\nfor (const p of candidates) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the two blank spacer paragraphs and the \"\u00a9 2020 ...\" footer\n# paragraph that used to immediately follow the\n# \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\" requirement line.\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOQ4086:*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOQ4086' anchor paragraph.\"\n}\n\n# The three paragraphs right after the anchor are: an empty spacer\n# paragraph, an empty page-break spacer paragraph, and the\n# \"\u00a9 2020 ... Contact: ... Jekyll ...\" copyright/footer paragraph.\n# Validate each before deleting so the script stays safe even if the\n# surrounding layout shifts slightly.\nfor ($k = 0; $k -lt 3; $k++) {\n    $p = $d.Paragraphs.Item($anchorIndex + 1)\n    $text = $p.Range.Text.Trim()\n    $isBlankSpacer = ($text -eq \"\")\n    $isCopyrightFooter = ($text -like \"*Contact:*\" -and $text -like \"*Jekyll*\")\n    if (-not ($isBlankSpacer -or $isCopyrightFooter)) {\n        throw \"Unexpected paragraph content while removing footer block: $text\"\n    }\n    $p.Range.Delete()\n}\n"}
